# fix(tokens), #113: Correct `sheet_id` definition.
#
# Adds a second worksheet literally named `!"` right after DATA, seeds its
# A1 with 0, points the selection at A2 on that sheet, updates DATA!B5's
# formula to reference '!"'!A1, and moves DATA's own selection to B6.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("DATA")

# Insert the new sheet right after DATA so it lands as the 2nd tab / sheetId 2.
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = '!"'
$newSheet.Range("A1").Value = 0
$newSheet.Range("A2").Select()

# Re-activate DATA (adding a sheet makes it the active one) and apply its edits.
$dataSheet.Activate()
$dataSheet.Range("B5").Formula = "=SUM(B2,B4,C2,'!""'!A1,IFERROR(BROKEN,0))"
$dataSheet.Range("B6").Select()
